$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 updates
$ws.Range("G8").Value = 2.88
$ws.Range("I8").Value = 2.45
$ws.Range("K8").Value = 2.05
$ws.Range("L8").Value = 3.2
$ws.Range("S8").Value = 2.08
$ws.Range("T8").Value = 1.73
$ws.Range("AA8").Value = 1.83
$ws.Range("AB8").Value = 1.83
$ws.Range("AC8").Value = 8.5
$ws.Range("AD8").Value = 13
$ws.Range("AJ8").Value = 6
$ws.Range("AO8").Value = 12
$ws.Range("AP8").Value = 10

# Row 10 updates
$ws.Range("G10").Value = 1.95
$ws.Range("H10").Value = 3.3
$ws.Range("I10").Value = 3.5
$ws.Range("J10").Value = 2.55
$ws.Range("K10").Value = 2.15
$ws.Range("L10").Value = 4
$ws.Range("M10").Value = 1.07
$ws.Range("N10").Value = 7
$ws.Range("O10").Value = 1.32
$ws.Range("P10").Value = 3.1
$ws.Range("S10").Value = 1.95
$ws.Range("T10").Value = 1.75
$ws.Range("W10").Value = 3.25
$ws.Range("X10").Value = 1.29
$ws.Range("Y10").Value = 1.39
$ws.Range("Z10").Value = 2.77
$ws.Range("AA10").Value = 1.8
$ws.Range("AB10").Value = 1.91
$ws.Range("AC10").Value = 7
$ws.Range("AF10").Value = 17.5
$ws.Range("AG10").Value = 16
$ws.Range("AH10").Value = 28
$ws.Range("AI10").Value = 7
$ws.Range("AJ10").Value = 6.4
$ws.Range("AL10").Value = 70
$ws.Range("AM10").Value = 600
$ws.Range("AN10").Value = 10
$ws.Range("AO10").Value = 19
$ws.Range("AS10").Value = 40

# Row 11 updates
$ws.Range("G11").Value = 2.2
$ws.Range("H11").Value = 3.4
$ws.Range("I11").Value = 3.2
$ws.Range("L11").Value = 3.75
$ws.Range("M11").Value = 1.05
$ws.Range("N11").Value = 11
$ws.Range("S11").Value = 1.9
$ws.Range("T11").Value = 1.95
$ws.Range("W11").Value = 3.25
$ws.Range("X11").Value = 1.33
$ws.Range("Y11").Value = 1.4
$ws.Range("Z11").Value = 2.75
$ws.Range("AA11").Value = 1.73
$ws.Range("AB11").Value = 2
$ws.Range("AC11").Value = 8.5
$ws.Range("AD11").Value = 11
$ws.Range("AH11").Value = 26
$ws.Range("AI11").Value = 11
$ws.Range("AM11").Value = 201
$ws.Range("AP11").Value = 12
$ws.Range("AS11").Value = 34
